# Add new "Drink" menu items (Sodas, Sweet Teas, Frappes, Iced Coffees) to
# the bottom of the MaxDonalds menu table on Sheet1, rows 49-59.
#
# NOTE on write order: the shared-string table is built in first-seen order.
# To reproduce the target sharedStrings.xml exactly, the *type* ("Drink")
# column of the first new row is written before its name column, and the
# "Medium Sweet Tea" label is written before "Small Sweet Tea" even though
# the latter ends up one row higher in the sheet. The numeric columns can be
# written in any order since they do not touch the shared-string table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Row 49: Small Soda ---------------------------------------------------
$ws.Cells.Item(49, 2).Value = "Drink"
$ws.Cells.Item(49, 1).Value = "Small Soda"
$ws.Cells.Item(49, 3).Value = 140
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 1.15

# --- Row 50: Medium Soda ---------------------------------------------------
$ws.Cells.Item(50, 1).Value = "Medium Soda"
$ws.Cells.Item(50, 2).Value = "Drink"
$ws.Cells.Item(50, 3).Value = 200
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 1.48

# --- Row 51: Large Soda ---------------------------------------------------
$ws.Cells.Item(51, 1).Value = "Large Soda"
$ws.Cells.Item(51, 2).Value = "Drink"
$ws.Cells.Item(51, 3).Value = 280
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 1.71

# --- Row 53 name first: Medium Sweet Tea (creates shared string before
#     "Small Sweet Tea" does, matching the author's original entry order) --
$ws.Cells.Item(53, 1).Value = "Medium Sweet Tea"

# --- Row 52: Small Sweet Tea ---------------------------------------------
$ws.Cells.Item(52, 1).Value = "Small Sweet Tea"
$ws.Cells.Item(52, 2).Value = "Drink"
$ws.Cells.Item(52, 3).Value = 150
$ws.Cells.Item(52, 4).Value = 0
$ws.Cells.Item(52, 5).Value = 0
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 10
$ws.Cells.Item(52, 8).Value = 1.15

# --- Row 53: Medium Sweet Tea (remaining columns) -------------------------
$ws.Cells.Item(53, 2).Value = "Drink"
$ws.Cells.Item(53, 3).Value = 180
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 10
$ws.Cells.Item(53, 8).Value = 1.48

# --- Row 54: Small Frappe (Type has a trailing space: "Drink ") -----------
$ws.Cells.Item(54, 1).Value = "Small Frappe"
$ws.Cells.Item(54, 2).Value = "Drink "
$ws.Cells.Item(54, 3).Value = 440
$ws.Cells.Item(54, 4).Value = 18
$ws.Cells.Item(54, 5).Value = 64
$ws.Cells.Item(54, 6).Value = 7
$ws.Cells.Item(54, 7).Value = 125
$ws.Cells.Item(54, 8).Value = 2.75

# --- Row 55: Medium Frappe -------------------------------------------------
$ws.Cells.Item(55, 1).Value = "Medium Frappe"
$ws.Cells.Item(55, 2).Value = "Drink"
$ws.Cells.Item(55, 3).Value = 540
$ws.Cells.Item(55, 4).Value = 22
$ws.Cells.Item(55, 5).Value = 79
$ws.Cells.Item(55, 6).Value = 9
$ws.Cells.Item(55, 7).Value = 160
$ws.Cells.Item(55, 8).Value = 3.32

# --- Row 56: Large Frappe ---------------------------------------------------
$ws.Cells.Item(56, 1).Value = "Large Frappe"
$ws.Cells.Item(56, 2).Value = "Drink"
$ws.Cells.Item(56, 3).Value = 670
$ws.Cells.Item(56, 4).Value = 26
$ws.Cells.Item(56, 5).Value = 97
$ws.Cells.Item(56, 6).Value = 11
$ws.Cells.Item(56, 7).Value = 190
$ws.Cells.Item(56, 8).Value = 3.9

# --- Row 57: Small Iced Coffee ---------------------------------------------
$ws.Cells.Item(57, 1).Value = "Small Iced Coffee"
$ws.Cells.Item(57, 2).Value = "Drink"
$ws.Cells.Item(57, 3).Value = 130
$ws.Cells.Item(57, 4).Value = 4.5
$ws.Cells.Item(57, 5).Value = 22
$ws.Cells.Item(57, 6).Value = 1
$ws.Cells.Item(57, 7).Value = 35
$ws.Cells.Item(57, 8).Value = 1.83

# --- Row 58: Medium Iced Coffee --------------------------------------------
$ws.Cells.Item(58, 1).Value = "Medium Iced Coffee"
$ws.Cells.Item(58, 2).Value = "Drink"
$ws.Cells.Item(58, 3).Value = 180
$ws.Cells.Item(58, 4).Value = 7
$ws.Cells.Item(58, 5).Value = 29
$ws.Cells.Item(58, 6).Value = 1
$ws.Cells.Item(58, 7).Value = 50
$ws.Cells.Item(58, 8).Value = 2.06

# --- Row 59: Large Iced Coffee ---------------------------------------------
$ws.Cells.Item(59, 1).Value = "Large Iced Coffee"
$ws.Cells.Item(59, 2).Value = "Drink"
$ws.Cells.Item(59, 3).Value = 260
$ws.Cells.Item(59, 4).Value = 9
$ws.Cells.Item(59, 5).Value = 43
$ws.Cells.Item(59, 6).Value = 2
$ws.Cells.Item(59, 7).Value = 65
$ws.Cells.Item(59, 8).Value = 2.29

# --- Reflect the view state that results from having just filled this
#     data in and landed on the next empty row -----------------------------
$ws.Range("E60").Select()
